# feat: add 2022-Q1 data
#
# Original tab order:  2020-Q4 | 2021-Q1 | 2021-Q2 | 2021-Q3 | 2021-Q4 | 总计
# New tab order:        2020-Q4 | 2021-Q1 | 2021-Q2 | 2021-Q3 | 2021-Q4 | 2022-Q1 | 总计
#
# The previous "总计" (totals/summary) sheet is repurposed into the new
# "2022-Q1" fund-holdings sheet (keeping its original sheetId), and a brand
# new "总计" sheet is appended right after it, holding the same rolling
# summary table as before plus a new 2022-Q1 row at the top (older rows
# shifted down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT, even when it looks like
# a number (e.g. "009411", "3.45") -- without leaving the cell's style index
# pointing at a "quotePrefix" variant. We do this by assigning with a leading
# apostrophe (forces text) and then pasting-formats-only from a pristine,
# never-touched blank cell to strip the quote-prefix style back to default.
# ---------------------------------------------------------------------------
$blankSrc = $wb.Worksheets.Item("2021-Q4").Range("ZZ1000")

function Set-TextCell($range, [string]$text) {
    $range.Value = "'" + $text
    $blankSrc.Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet into the new "2022-Q1" sheet, but
# first read out its existing summary rows so we can reuse them later.
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")

$oldDate1 = $oldTotal.Range("B2").Value2
$oldCount1 = $oldTotal.Range("C2").Value2
$oldValue1 = $oldTotal.Range("D2").Value2

$oldDate2 = $oldTotal.Range("B3").Value2
$oldCount2 = $oldTotal.Range("C3").Value2
$oldValue2 = $oldTotal.Range("D3").Value2

$oldDate3 = $oldTotal.Range("B4").Value2
$oldCount3 = $oldTotal.Range("C4").Value2
$oldValue3 = $oldTotal.Range("D4").Value2

$oldDate4 = $oldTotal.Range("B5").Value2
$oldCount4 = $oldTotal.Range("C5").Value2
$oldValue4 = $oldTotal.Range("D5").Value2

$oldDate5 = $oldTotal.Range("B6").Value2
$oldCount5 = $oldTotal.Range("C6").Value2
$oldValue5 = $oldTotal.Range("D6").Value2

$q1 = $oldTotal
$q1.Name = "2022-Q1"

# Clear out the old summary-table contents (A1:D6) before writing new data.
$q1.Range("A1:D6").Clear()

# Header row (matches the other quarter sheets).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy header style/formatting from an existing quarter sheet so the new
# header row matches the look of the others (bold, border, centered).
$wb.Worksheets.Item("2021-Q4").Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Row 2: 009411 中银科技创新一年定期开放混合
$q1.Range("A2").Value2 = 0
Set-TextCell $q1.Range("B2") "009411"
Set-TextCell $q1.Range("C2") "中银科技创新一年定期开放混合"
Set-TextCell $q1.Range("D2") "3.45"
Set-TextCell $q1.Range("E2") "90.71"
Set-TextCell $q1.Range("F2") "3.31"
Set-TextCell $q1.Range("G2") "0.1142"
$q1.Range("H2").Value2 = 9

# Row 3: 163809 中银蓝筹精选灵活配置混合
$q1.Range("A3").Value2 = 1
Set-TextCell $q1.Range("B3") "163809"
Set-TextCell $q1.Range("C3") "中银蓝筹精选灵活配置混合"
Set-TextCell $q1.Range("D3") "3.36"
Set-TextCell $q1.Range("E3") "79.31"
Set-TextCell $q1.Range("F3") "2.77"
Set-TextCell $q1.Range("G3") "0.0931"
$q1.Range("H3").Value2 = 9

# Style the A column index cells (row number style, same as other sheets).
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)
$q1.Range("A2").Value2 = 0
$q1.Range("A3").Value2 = 1

# ---------------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet after "2022-Q1" holding the updated
# rolling summary table (2022-Q1 added at the top, older rows shifted down).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Copy the header style from the original summary sheet look.
$wb.Worksheets.Item("2021-Q4").Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value2 = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value2 = 2
$total.Range("D2").Value2 = 0.21

$total.Range("A3").Value2 = 1
$total.Range("B3").Value = $oldDate1
$total.Range("C3").Value2 = $oldCount1
$total.Range("D3").Value2 = $oldValue1

$total.Range("A4").Value2 = 2
$total.Range("B4").Value = $oldDate2
$total.Range("C4").Value2 = $oldCount2
$total.Range("D4").Value2 = $oldValue2

$total.Range("A5").Value2 = 3
$total.Range("B5").Value = $oldDate3
$total.Range("C5").Value2 = $oldCount3
$total.Range("D5").Value2 = $oldValue3

$total.Range("A6").Value2 = 4
$total.Range("B6").Value = $oldDate4
$total.Range("C6").Value2 = $oldCount4
$total.Range("D6").Value2 = $oldValue4

$total.Range("A7").Value2 = 5
$total.Range("B7").Value = $oldDate5
$total.Range("C7").Value2 = $oldCount5
$total.Range("D7").Value2 = $oldValue5

# Style the A column index cells (0..5) to match the other summary sheets.
$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$total.Range("A2").Value2 = 0
$total.Range("A3").Value2 = 1
$total.Range("A4").Value2 = 2
$total.Range("A5").Value2 = 3
$total.Range("A6").Value2 = 4
$total.Range("A7").Value2 = 5
